$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deskcount")

# "Include in Occupancy Calculation" (column F): Yes -> No
#   row 17 = Greenwood Village, row 39 = Tampa, row 48 = Santiago, row 49 = Sao Paulo
$ws.Range("F17").Value = "No"
$ws.Range("F39").Value = "No"
$ws.Range("F48").Value = "No"
$ws.Range("F49").Value = "No"

# Deskcount (column C) for Melbourne (row 45): 30 -> 32
$ws.Range("C45").Value = 32

# Reflect the saved view state (scrolled down, C46 selected)
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C46").Select()
